$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.47%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'36.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.46%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.023"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.30%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07880"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.49%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.140"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-3.17%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "'GateToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'4.132"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.60%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'KuCoinToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'7.940"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.72%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9222"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.67%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09738"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.20%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.03%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08592"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.29%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03578"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.85%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09923"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.26%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001441"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-3.84%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.005692"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.70%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.471"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.26%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'BTSEToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'2.752"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'19.05%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'BitpandaEcosystemToken"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'0.3375"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.73%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'ProBitToken"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.1347"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.52%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'MCDex"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'5.146"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'7.39%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'ZBToken"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.2208"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.30%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'CoinExToken"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.04571"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.40%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-1.61%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004807"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-7.58%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-7.17%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004752"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'74.87%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01853"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.91%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04724"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.84%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007799"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.63%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1387"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.06%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007719"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.08%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-3.93%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01139"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'9.43%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006373"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.22%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'0.19%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'52.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'46.76%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.001901"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-29.31%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.04%"
$ws.Range("E51").Style = "Normal"
